# Nexial "#system" showcase sheet update:
#   - [base commands] add NEW command `assertMatch(text,regex)` (inserted
#     alphabetically into the "base" list, column F)
#   - [external commands] add NEW command `openFile(filePath)` (inserted
#     alphabetically into the "external" list, column J)
#   - remove the obsolete "tn.5250" column (column AA) from the lookup
#     sheet, shifting the web/webalert/webcookie/ws/ws.async/xml lists one
#     column to the left (AB->AA, AC->AB, AD->AC, AE->AD, AF->AE, AG->AF)
#   - the "target" category list (column A) drops the "tn.5250" entry

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) "base" commands (column F) - insert assertMatch(text,regex) so the
#    whole list stays alphabetically sorted; write the full resulting
#    list back out (F2:F45).
# ---------------------------------------------------------------------
$baseCommands = @(
    "appendText(var,appendWith)",
    "assertArrayContain(array,expected)",
    "assertArrayEqual(array1,array2,exactOrder)",
    "assertArrayNotContain(array,unexpected)",
    "assertContains(text,substring)",
    "assertCount(text,regex,expects)",
    "assertEmpty(text)",
    "assertEndsWith(text,suffix)",
    "assertEqual(expected,actual)",
    "assertMatch(text,regex)",
    "assertNotContain(text,substring)",
    "assertNotEmpty(text)",
    "assertNotEqual(expected,actual)",
    "assertStartsWith(text,prefix)",
    "assertTextOrder(var,descending)",
    "assertVarNotPresent(var)",
    "assertVarPresent(var)",
    "clear(vars)",
    "clearClipboard()",
    "copyFromClipboard(var)",
    "copyIntoClipboard(text)",
    "failImmediate(text)",
    "incrementChar(var,amount,config)",
    "macro(file,sheet,name)",
    "macroFlex(macro,input,output)",
    "outputToCloud(resource)",
    "prependText(var,prependWith)",
    "repeatUntil(steps,maxWaitMs)",
    "save(var,value)",
    "saveCount(text,regex,saveVar)",
    "saveMatches(text,regex,saveVar)",
    "saveReplace(text,regex,replace,saveVar)",
    "saveVariablesByPrefix(var,prefix)",
    "saveVariablesByRegex(var,regex)",
    "section(steps)",
    "split(text,delim,saveVar)",
    "startRecording()",
    "stopRecording()",
    "substringAfter(text,delim,saveVar)",
    "substringBefore(text,delim,saveVar)",
    "substringBetween(text,start,end,saveVar)",
    "verbose(text)",
    "waitFor(waitMs)",
    "waitForCondition(conditions,maxWaitMs)"
)
for ($i = 0; $i -lt $baseCommands.Count; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $baseCommands[$i]
}

# ---------------------------------------------------------------------
# 2) "external" commands (column J) - insert openFile(filePath); write
#    the full resulting list back out (J2:J7).
# ---------------------------------------------------------------------
$externalCommands = @(
    "openFile(filePath)",
    "runJUnit(className)",
    "runProgram(programPathAndParams)",
    "runProgramNoWait(programPathAndParams)",
    "tail(id,file)",
    "terminate(programName)"
)
for ($i = 0; $i -lt $externalCommands.Count; $i++) {
    $ws.Cells.Item($i + 2, 10).Value = $externalCommands[$i]
}

# ---------------------------------------------------------------------
# 3) Drop the "tn.5250" column (AA) entirely - everything to its right
#    (web, webalert, webcookie, ws, ws.async, xml) shifts one column left.
# ---------------------------------------------------------------------
$ws.Columns.Item(27).Delete()

# ---------------------------------------------------------------------
# 4) "target" category list (column A) loses the "tn.5250" row.
# ---------------------------------------------------------------------
$ws.Range("A27").Delete(-4162)

# ---------------------------------------------------------------------
# 5) Fix up the named ranges to reflect the new sizes/positions.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
